$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Branch" column before column L (12th column)
$ws.Range("L1").EntireColumn.Insert()

# Populate the new column's header and data
$ws.Range("L1").Value = "Branch"
$ws.Range("L2").Value = "CHOWGULE LAVGAN"
$ws.Range("L3").Value = "CHOWGULE LAVGAN"

# The AutoFilter's hidden _FilterDatabase defined name needs to extend
# to the new last column (AQ instead of AP)
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$AQ`$1"

# Match the new selection location
$ws.Range("L1").Select()
